# Stock_Analysis.xlsx edit
# Swap the "Stock Name" (A) and "Ticker" (B) columns on Potential_Investments,
# since the Ticker column is what the external scraper iterates over and
# therefore needs to be first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Potential_Investments")

# ------------------------------------------------------------------
# 1. Swap columns A and B (entire columns - content, formatting, width)
#    Select column A, cut it, and re-insert it after column B - this is
#    the "drag column with Shift held" move done via COM.
# ------------------------------------------------------------------
$ws.Columns.Item(1).Cut() | Out-Null
$ws.Columns.Item(3).Insert() | Out-Null

# ------------------------------------------------------------------
# 2. The sheet's AutoFilter / _FilterDatabase used to span A1:N18
#    (Stock Name .. Growth). Now that the Ticker column moved to A,
#    the filtered/sorted range should start at column B instead.
#    Re-create the AutoFilter over B1:N18.
#
#    Row 19 sits directly below the filtered block with no gap, so a
#    freshly-applied AutoFilter would auto-expand to include it. Stash
#    row 19's contents, blank it out for the re-apply, then restore it
#    so the final data is untouched.
# ------------------------------------------------------------------
$lastRow = $ws.Range("B1:N18").Rows.Count + 1
$stashRange = $ws.Range("A19:P19")
$stashedRow19 = $stashRange.Formula

$ws.AutoFilterMode = $false
$stashRange.ClearContents() | Out-Null

$ws.Range("B1:N18").AutoFilter() | Out-Null

$stashRange.Formula = $stashedRow19

# ------------------------------------------------------------------
# 3. Keep the hidden _xlnm._FilterDatabase defined name lined up with
#    the AutoFilter's new range.
# ------------------------------------------------------------------
$wb.Names.Item("Potential_Investments!_FilterDatabase").RefersTo = "=Potential_Investments!`$B`$1:`$N`$18"

# ------------------------------------------------------------------
# 4. Leave the sheet with column A selected (whole-sheet selection, as
#    happens right after dragging a column into a new position).
# ------------------------------------------------------------------
$ws.Cells.Select() | Out-Null
